# "Adding properties to folders"
# Populate the Value column (B) for the folder-path related settings on the
# "Constants" sheet and add a new "URL_Test" entry, then tidy up the view
# state (column autosize, zoom, selection) to match the authored session.

$wb = $excel.ActiveWorkbook

$settings  = $wb.Worksheets.Item("Settings")
$constants = $wb.Worksheets.Item("Constants")

# Note: the "Assets" sheet is intentionally left untouched - it only had
# cosmetic/view changes in the original session.

# --- Constants sheet: add the new URL_Test row and fill in folder values ---
# (entered in the same order the author did: new row first, then the
# existing rows bottom-to-top)
$constants.Range("A7").Value = "URL_Test"
$constants.Range("B7").Value = "https://forms.office.com/Pages/ResponsePage.aspx?id=x8fjOlHq3kaVQBL0EQ6smyFnQk63wvpEhjCaoYxMrW5UQlpSVFVGSlpRWE8xQzFVTEVBMzdZVzRaTC4u"

$constants.Range("B6").Value = "C:\Users\ivan_\OneDrive\Documentos\Power Automate\Repo\Proyecto\Junior - UNO"
$constants.Range("B5").Value = "C:\Users\ivan_\OneDrive\Documentos\Power Automate\Repo\Proyecto\Junior - UNO\Output\LogFile"
$constants.Range("B4").Value = "C:\Users\ivan_\OneDrive\Documentos\Power Automate\Repo\Proyecto\Junior - UNO\Output\ScreenshotError"

# Column B now holds long path/URL values - autosize it like the author did
$constants.Columns.Item(2).AutoFit() | Out-Null

# Move the cursor to where the author left it
$constants.Range("B8").Select() | Out-Null

# --- Settings sheet: zoom was reduced for the new window size ---
$settings.Activate()
$excel.ActiveWindow.Zoom = 115

# Return focus to Constants, which stays the active tab
$constants.Activate()

Write-Host "Applied 'Adding properties to folders' changes"
